$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original (pre-edit) values for the columns that rotate/change
# (D: Fecha, J: Volumen, K: Precio minimo, L: Precio maximo, M: Precio promedio,
#  O: Origen, P: Precio $/Kg), keyed by the original source row, then write them
# into their new target rows per the edit.
$src = @{}
$src[2] = @{ "D"=44335; "J"=100; "K"=18000; "L"=20000; "M"=19000; "O"="Provincia de Limarí"; "P"=760 }
$src[3] = @{ "D"=44497; "J"=150; "K"=6000; "L"=6500; "M"=6333; "O"="Región Metropolitana"; "P"=253 }
$src[4] = @{ "D"=44188; "J"=100; "K"=18000; "L"=20000; "M"=19000; "O"="Región Metropolitana"; "P"=760 }
$src[5] = @{ "D"=44537; "J"=160; "K"=8500; "L"=9000; "M"=8719; "O"="Región del Maule"; "P"=349 }
$src[6] = @{ "D"=44476; "J"=100; "K"=7000; "L"=7500; "M"=7250; "O"="Región Metropolitana"; "P"=290 }
$src[7] = @{ "D"=44503; "J"=250; "K"=9000; "L"=10000; "M"=9400; "O"="Provincia de Melipilla"; "P"=376 }
$src[8] = @{ "D"=44509; "J"=100; "K"=6500; "L"=7000; "M"=6750; "O"="Región Metropolitana"; "P"=270 }
$src[9] = @{ "D"=44483; "J"=350; "K"=5500; "L"=6000; "M"=5714; "O"="Región Metropolitana"; "P"=229 }
$src[10] = @{ "D"=44498; "J"=220; "K"=7000; "L"=7500; "M"=7273; "O"="Región Metropolitana"; "P"=291 }
$src[11] = @{ "D"=44523; "J"=100; "K"=9000; "L"=10000; "M"=9500; "O"="Región Metropolitana"; "P"=380 }
$src[12] = @{ "D"=44545; "J"=140; "K"=14000; "L"=15000; "M"=14429; "O"="Provincia de Chacabuco"; "P"=577 }
$src[13] = @{ "D"=44384; "J"=100; "K"=12000; "L"=13000; "M"=12500; "O"="Región de Coquimbo"; "P"=500 }
$src[14] = @{ "D"=44526; "J"=100; "K"=7500; "L"=8000; "M"=7750; "O"="Región Metropolitana"; "P"=310 }
$src[15] = @{ "D"=44467; "J"=100; "K"=8000; "L"=9000; "M"=8500; "O"="Región Metropolitana"; "P"=340 }
$src[16] = @{ "D"=44692; "J"=100; "K"=20000; "L"=22000; "M"=21000; "O"="Región Metropolitana"; "P"=840 }
$src[17] = @{ "D"=44533; "J"=180; "K"=8000; "L"=8500; "M"=8222; "O"="Región del Maule"; "P"=329 }
$src[18] = @{ "D"=44517; "J"=130; "K"=6000; "L"=6500; "M"=6269; "O"="Región Metropolitana"; "P"=251 }
$src[19] = @{ "D"=44461; "J"=100; "K"=13000; "L"=14000; "M"=13500; "O"="Provincia del Elquí"; "P"=540 }
$src[20] = @{ "D"=44160; "J"=100; "K"=9000; "L"=10000; "M"=9500; "O"="Región Metropolitana"; "P"=380 }
$src[21] = @{ "D"=44162; "J"=100; "K"=7500; "L"=8000; "M"=7750; "O"="Región Metropolitana"; "P"=310 }
$src[22] = @{ "D"=44540; "J"=140; "K"=11000; "L"=12000; "M"=11429; "O"="Región del Maule"; "P"=457 }
$src[23] = @{ "D"=44482; "J"=430; "K"=8000; "L"=8500; "M"=8267; "O"="Región de O'Higgins"; "P"=331 }
$src[24] = @{ "D"=44316; "J"=100; "K"=16000; "L"=18000; "M"=17000; "O"="Región Metropolitana"; "P"=680 }
$src[25] = @{ "D"=44505; "J"=180; "K"=6000; "L"=6500; "M"=6222; "O"="Región del Maule"; "P"=249 }
$src[26] = @{ "D"=44673; "J"=80; "K"=18000; "L"=19000; "M"=18375; "O"="Región Metropolitana"; "P"=735 }
$src[27] = @{ "D"=44351; "J"=100; "K"=15000; "L"=16000; "M"=15500; "O"="Región Metropolitana"; "P"=620 }
$src[28] = @{ "D"=44454; "J"=100; "K"=13000; "L"=14000; "M"=13500; "O"="Provincia del Elquí"; "P"=540 }

# Target-row -> source-row mapping
$map = @{}
$map[2] = 5
$map[3] = 8
$map[4] = 23
$map[5] = 4
$map[6] = 22
$map[7] = 19
$map[8] = 15
$map[9] = 6
$map[10] = 14
$map[11] = 12
$map[12] = 18
$map[13] = 21
$map[14] = 9
$map[15] = 3
$map[16] = 26
$map[17] = 7
$map[18] = 2
$map[19] = 13
$map[20] = 28
$map[21] = 10
$map[22] = 24
$map[23] = 20
$map[24] = 25
$map[25] = 17
$map[26] = 11
$map[28] = 16

foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $vals = $src[$oldRow]
    $ws.Range("D" + $newRow).Value = $vals["D"]
    $ws.Range("J" + $newRow).Value = $vals["J"]
    $ws.Range("K" + $newRow).Value = $vals["K"]
    $ws.Range("L" + $newRow).Value = $vals["L"]
    $ws.Range("M" + $newRow).Value = $vals["M"]
    $ws.Range("O" + $newRow).Value = $vals["O"]
    $ws.Range("P" + $newRow).Value = $vals["P"]
}